$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear out the header cells that are no longer used (TGL_LAHIR..NAMA_PEGAWAI)
$ws.Range("E1:Q1").ClearContents()

# Clear out the row 2 data cells that are no longer used, but keep D2/E2/N2
# (they stay as empty, styled cells)
$ws.Range("E2:Q2").ClearContents()

# New row 2 values
$ws.Range("A2").Value = 99
$ws.Range("B2").Value = "yuyuk"
